# Swap the taxon-observation data between row 4 and row 6, while leaving
# the shared/location columns (C, P, Q, R, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AT, AW, AX, AY) untouched.
#
# Row 4 currently holds the "Buskskvätta" (bird) record and Row 6 holds the
# "Hällebräcka" (plant) record; after the edit they trade places. In
# addition, the "Publik kommentar" (AC) and "Fritext" (AI) comments that
# belong to the Hällebräcka record move from row 6 to row 4 along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture row 4's original values (the Buskskvätta record) ---
$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()
$d4 = $ws.Range("D4").Value()
$e4 = $ws.Range("E4").Value()
$f4 = $ws.Range("F4").Value()
$g4 = $ws.Range("G4").Value()
$h4 = $ws.Range("H4").Value()
$s4 = $ws.Range("S4").Value()

# --- Capture row 6's original values (the Hällebräcka record) ---
$a6 = $ws.Range("A6").Value()
$b6 = $ws.Range("B6").Value()
$d6 = $ws.Range("D6").Value()
$e6 = $ws.Range("E6").Value()
$f6 = $ws.Range("F6").Value()
$g6 = $ws.Range("G6").Value()
$h6 = $ws.Range("H6").Value()
$s6 = $ws.Range("S6").Value()
$ac6 = $ws.Range("AC6").Value()
$ai6 = $ws.Range("AI6").Value()

# --- Write row 6's original values into row 4 ---
$ws.Range("A4").Value = $a6
$ws.Range("B4").Value = $b6
$ws.Range("D4").Value = $d6
$ws.Range("E4").Value = $e6
$ws.Range("F4").Value = $f6
$ws.Range("G4").Value = $g6
$ws.Range("H4").Value = $h6
$ws.Range("S4").Value = $s6
$ws.Range("AC4").Value = $ac6
$ws.Range("AI4").Value = $ai6

# --- Write row 4's original values into row 6 ---
$ws.Range("A6").Value = $a4
$ws.Range("B6").Value = $b4
$ws.Range("D6").Value = $d4
$ws.Range("E6").Value = $e4
$ws.Range("F6").Value = $f4
$ws.Range("G6").Value = $g4
$ws.Range("H6").Value = $h4
$ws.Range("S6").Value = $s4

# Row 6 no longer carries the comment/fritext cells (they moved to row 4).
$ws.Range("AC6").Value = ""
$ws.Range("AI6").Value = ""
